$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2593.6155
$ws.Range("I40").Value = 2303
$ws.Range("J40").Value = 2842.7144
$ws.Range("K40").Value = 2303
$ws.Range("L40").Value = 2842.7144
$ws.Range("M40").Value = -2128
$ws.Range("N40").Value = -3192.7144
$ws.Range("H64").Value = 9653.714
$ws.Range("I64").Value = 7687
$ws.Range("J64").Value = 9981.5
$ws.Range("K64").Value = 7687
$ws.Range("L64").Value = 9981.5
$ws.Range("M64").Value = -7439
$ws.Range("N64").Value = -10477.5
$ws.Range("H67").Value = 9653.714
$ws.Range("I67").Value = 7687
$ws.Range("J67").Value = 9981.5
$ws.Range("K67").Value = 7687
$ws.Range("L67").Value = 9981.5
$ws.Range("M67").Value = -6829
$ws.Range("N67").Value = -11697.5
$ws.Range("H74").Value = 5864.8965
$ws.Range("I74").Value = 4917.875
$ws.Range("J74").Value = 7030.4614
$ws.Range("K74").Value = 4917.875
$ws.Range("L74").Value = 7030.4614
$ws.Range("M74").Value = -3981.875
$ws.Range("N74").Value = -8902.4614
$ws.Range("H77").Value = 5864.8965
$ws.Range("I77").Value = 4917.875
$ws.Range("J77").Value = 7030.4614
$ws.Range("K77").Value = 24589.375
$ws.Range("L77").Value = 35152.307
$ws.Range("M77").Value = -19909.375
$ws.Range("N77").Value = -44512.307
$ws.Range("H98").Value = 8402.200000000001
$ws.Range("I98").Value = 622.1905
$ws.Range("J98").Value = 49247.25
$ws.Range("K98").Value = 622.1905
$ws.Range("L98").Value = 49247.25
$ws.Range("M98").Value = 875.8095
$ws.Range("N98").Value = -52243.25
$ws.Range("H122").Value = 8402.200000000001
$ws.Range("I122").Value = 622.1905
$ws.Range("J122").Value = 49247.25
$ws.Range("K122").Value = 1866.5715
$ws.Range("L122").Value = 147741.75
$ws.Range("M122").Value = 583.4285
$ws.Range("N122").Value = -152641.75
$ws.Range("H137").Value = 2658.2856
$ws.Range("I137").Value = 1967.25
$ws.Range("J137").Value = 3083.5386
$ws.Range("K137").Value = 5901.75
$ws.Range("L137").Value = 9250.6158
$ws.Range("M137").Value = -3351.75
$ws.Range("N137").Value = -14350.6158
$ws.Range("H138").Value = 4320.302
$ws.Range("I138").Value = 3566.652
$ws.Range("K138").Value = 10699.956
$ws.Range("M138").Value = -5559.956

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1779.6216
$ws.Range("I32").Value = 1016.5507
$ws.Range("K32").Value = 1016.5507
$ws.Range("M32").Value = -729.5507
$ws.Range("H61").Value = 2040.0834
$ws.Range("I61").Value = 2180.1365
$ws.Range("J61").Value = 499.5
$ws.Range("K61").Value = 2180.1365
$ws.Range("L61").Value = 499.5
$ws.Range("M61").Value = -1968.1365
$ws.Range("N61").Value = -923.5
$ws.Range("H63").Value = 3136
$ws.Range("I63").Value = 3236.4
$ws.Range("J63").Value = 2968.6667
$ws.Range("K63").Value = 3236.4
$ws.Range("L63").Value = 2968.6667
$ws.Range("M63").Value = -2550.4
$ws.Range("N63").Value = -4340.6667
$ws.Range("H66").Value = 3136
$ws.Range("I66").Value = 3236.4
$ws.Range("J66").Value = 2968.6667
$ws.Range("K66").Value = 16182
$ws.Range("L66").Value = 14843.3335
$ws.Range("M66").Value = -12750
$ws.Range("N66").Value = -21707.3335
$ws.Range("H74").Value = 2346.853
$ws.Range("I74").Value = 1687.625
$ws.Range("K74").Value = 1687.625
$ws.Range("M74").Value = -813.625
$ws.Range("H77").Value = 2346.853
$ws.Range("I77").Value = 1687.625
$ws.Range("K77").Value = 8438.125
$ws.Range("M77").Value = -4070.125
$ws.Range("H132").Value = 3499.5652
$ws.Range("I132").Value = 2918.389
$ws.Range("K132").Value = 8755.167000000001
$ws.Range("M132").Value = -6225.167000000001
$ws.Range("H136").Value = 2040.0834
$ws.Range("I136").Value = 2180.1365
$ws.Range("J136").Value = 499.5
$ws.Range("K136").Value = 6540.4095
$ws.Range("L136").Value = 1498.5
$ws.Range("M136").Value = -3990.4095
$ws.Range("N136").Value = -6598.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3391
$ws.Range("I105").Value = 3267.2307
$ws.Range("J105").Value = 5000
$ws.Range("K105").Value = 3267.2307
$ws.Range("L105").Value = 5000
$ws.Range("M105").Value = -1520.2307
$ws.Range("N105").Value = -8494
$ws.Range("H134").Value = 4976.259
$ws.Range("I134").Value = 4182.7676
$ws.Range("J134").Value = 7250.933
$ws.Range("K134").Value = 12548.3028
$ws.Range("L134").Value = 21752.799
$ws.Range("M134").Value = -10013.3028
$ws.Range("N134").Value = -26822.799

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1532.4054
$ws.Range("I58").Value = 1063
$ws.Range("J58").Value = 2399
$ws.Range("K58").Value = 1063
$ws.Range("L58").Value = 2399
$ws.Range("M58").Value = -860
$ws.Range("N58").Value = -2805
$ws.Range("H59").Value = 71110.89
$ws.Range("J59").Value = 71110.89
$ws.Range("L59").Value = 71110.89
$ws.Range("N59").Value = -73400.89
$ws.Range("H62").Value = 4248.8
$ws.Range("I62").Value = 3480
$ws.Range("J62").Value = 5017.6
$ws.Range("K62").Value = 3480
$ws.Range("L62").Value = 5017.6
$ws.Range("M62").Value = -2856
$ws.Range("N62").Value = -6265.6
$ws.Range("H65").Value = 4248.8
$ws.Range("I65").Value = 3480
$ws.Range("J65").Value = 5017.6
$ws.Range("K65").Value = 17400
$ws.Range("L65").Value = 25088
$ws.Range("M65").Value = -14280
$ws.Range("N65").Value = -31328
$ws.Range("H134").Value = 3693.0605
$ws.Range("I134").Value = 3062.4666
$ws.Range("K134").Value = 9187.399800000001
$ws.Range("M134").Value = -6652.399800000001
$ws.Range("H136").Value = 1532.4054
$ws.Range("I136").Value = 1063
$ws.Range("J136").Value = 2399
$ws.Range("K136").Value = 3189
$ws.Range("L136").Value = 7197
$ws.Range("M136").Value = -639
$ws.Range("N136").Value = -12297

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 3001.6
$ws.Range("I132").Value = 2004
$ws.Range("J132").Value = 3251
$ws.Range("K132").Value = 18036
$ws.Range("L132").Value = 29259
$ws.Range("M132").Value = -15506
$ws.Range("N132").Value = -34319

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7921.7856
$ws.Range("I70").Value = 7616.684
$ws.Range("J70").Value = 8565.888999999999
$ws.Range("K70").Value = 7616.684
$ws.Range("L70").Value = 8565.888999999999
$ws.Range("M70").Value = -7346.684
$ws.Range("N70").Value = -9105.888999999999
$ws.Range("H73").Value = 7921.7856
$ws.Range("I73").Value = 7616.684
$ws.Range("J73").Value = 8565.888999999999
$ws.Range("K73").Value = 7616.684
$ws.Range("L73").Value = 8565.888999999999
$ws.Range("M73").Value = -6680.684
$ws.Range("N73").Value = -10437.889
$ws.Range("H132").Value = 3406.9285
$ws.Range("I132").Value = 3634.2222
$ws.Range("J132").Value = 2997.8
$ws.Range("K132").Value = 10902.6666
$ws.Range("L132").Value = 8993.400000000001
$ws.Range("M132").Value = -8372.6666
$ws.Range("N132").Value = -14053.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1586.3572
$ws.Range("I61").Value = 1700.6923
$ws.Range("J61").Value = 100
$ws.Range("K61").Value = 1700.6923
$ws.Range("L61").Value = 100
$ws.Range("M61").Value = -1498.6923
$ws.Range("N61").Value = -504
$ws.Range("H113").Value = 1586.3572
$ws.Range("I113").Value = 1700.6923
$ws.Range("J113").Value = 100
$ws.Range("K113").Value = 1700.6923
$ws.Range("L113").Value = 100
$ws.Range("M113").Value = 469.3077000000001
$ws.Range("N113").Value = -4440

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3084.7637
$ws.Range("I132").Value = 3050.0208
$ws.Range("J132").Value = 3323
$ws.Range("K132").Value = 9150.062399999999
$ws.Range("L132").Value = 9969
$ws.Range("M132").Value = -6620.062399999999
$ws.Range("N132").Value = -15029
$ws.Range("H136").Value = 2544.8276
$ws.Range("J136").Value = 2111.7144
$ws.Range("L136").Value = 6335.1432
$ws.Range("N136").Value = -11435.1432
